$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15: update First_Detection_Image
$ws.Range("D15").Value = "image_20250807110238_ppp0.jpg"

# Row 16: update First_Detection_Image, First_Coords and First_Confidence
$ws.Range("D16").Value = "image_20250808100711_ppp0.jpg"
$ws.Range("I16").Value = "1182,409,1232,451"

# J16 is a numeric-looking text value ("0.75"); force it to stay text
# (matching the original inlineStr cell type) and then clear the
# quote-prefix style that forcing text introduces, so the cell keeps
# the default (unstyled) formatting like the rest of the column.
$ws.Range("J16").Value = "'0.75"
$ws.Range("J16").Style = "Normal"
